# Generate Report for Handoff
# Updates the "Latest Handoff/HO Xliff Generate" timestamps and the
# "Priority" column (marks several rows as hand-off type "ht") for the
# rows that were re-handed-off: 8, 9, 11, 12, 13, 14.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 11, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-01 22:23:20"
}

# --- zh-cn sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-01 22:23:15"
}

# --- de-de sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-01 22:23:20"
}
